$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values are plain numeric-looking strings (single decimal point)
# must be forced to Text format first, otherwise Excel auto-converts them to
# actual numbers (losing trailing zeros / exact formatting) on assignment.

$ws.Range("D2").Value = "35.310.07"
$ws.Range("E2").Value = "  -0.15%  "
$ws.Range("D3").Value = "1.878.51"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  -0.62%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.64"
$ws.Range("E5").Value = "  -3.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.687"
$ws.Range("E6").Value = "  -5.88%  "
$ws.Range("E7").Value = "  -0.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "43.14"
$ws.Range("E8").Value = "  +5.58%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.350"
$ws.Range("E9").Value = "  -5.63%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0736"
$ws.Range("E10").Value = "  -3.11%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0970"
$ws.Range("E11").Value = "  -1.89%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "13.09"
$ws.Range("E12").Value = "  +1.27%  "
$ws.Range("D13").Value = "2.152.42"
$ws.Range("E13").Value = "  -1.49%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.733"
$ws.Range("E14").Value = "  +0.16%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.94"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("D16").Value = "1.885.87"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").Value = "35.356.96"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "73.44"
$ws.Range("E18").Value = "  -2.46%  "
$ws.Range("D19").Value = "0.0₃0820"
$ws.Range("E19").Value = "  -3.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.53"
$ws.Range("E20").Value = "  +0.19%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.79"
$ws.Range("E21").Value = "  -2.48%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.93"
$ws.Range("E22").Value = "  -3.86%  "
$ws.Range("E23").Value = "  -0.78%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.56"
$ws.Range("E24").Value = "  +4.95%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.18"
$ws.Range("E25").Value = "  -11.27%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "165.28"
$ws.Range("E26").Value = "  -0.57%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.44"
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "18.27"
$ws.Range("E28").Value = "  -3.02%  "
$ws.Range("E29").Value = "  -4.67%  "
$ws.Range("E30").Value = "  -0.01%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.73"
$ws.Range("E31").Value = "  +4.73%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.23"
$ws.Range("E32").Value = "  -3.41%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0578"
$ws.Range("E33").Value = "  -1.90%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.20"
$ws.Range("E34").Value = "  -1.70%  "
$ws.Range("E35").Value = "  -0.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.849"
$ws.Range("E36").Value = "  -7.30%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.98"
$ws.Range("E37").Value = "  -2.80%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.58"
$ws.Range("E38").Value = "  -20.74%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0690"
$ws.Range("E39").Value = "  +6.49%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "17.04"
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "97.09"
$ws.Range("E41").Value = "  -0.20%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.0214"
$ws.Range("E42").Value = "  -3.26%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.08"
$ws.Range("E43").Value = "  -4.20%  "
$ws.Range("D44").Value = "1.286.47"
$ws.Range("E44").Value = "  -3.81%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.33"
$ws.Range("E45").Value = "  -5.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0807"
$ws.Range("E46").Value = "  +6.78%  "
$ws.Range("E47").Value = "  -1.11%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.73"
$ws.Range("E48").Value = "  -1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "12.08"
$ws.Range("E49").Value = "  +2.81%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "43.12"
$ws.Range("E50").Value = "  -4.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.27"
$ws.Range("E51").Value = "  -6.89%  "
